$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 20
$ws.Cells.Item(20, 2).Value = 6799253
$ws.Cells.Item(20, 5).Value = 'Hillerd'
$ws.Cells.Item(20, 6).Value = 'Hobro IK'
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 1
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 'A'
$ws.Cells.Item(20, 12).Value = 2.6
$ws.Cells.Item(20, 13).Value = 3.5
$ws.Cells.Item(20, 14).Value = 2.5
$ws.Cells.Item(20, 15).Value = 3
$ws.Cells.Item(20, 16).Value = 3.6
$ws.Cells.Item(20, 17).Value = 2.25
$ws.Cells.Item(20, 18).Value = 0.25
$ws.Cells.Item(20, 19).Value = 1.925
$ws.Cells.Item(20, 20).Value = 1.925
$ws.Cells.Item(20, 21).Value = 2.5
$ws.Cells.Item(20, 22).Value = 1.825
$ws.Cells.Item(20, 23).Value = 2.025
$ws.Cells.Item(20, 24).Value = -1
$ws.Cells.Item(20, 25).Value = -1
$ws.Cells.Item(20, 26).Value = 1.25
$ws.Cells.Item(20, 27).Value = -1
$ws.Cells.Item(20, 28).Value = 0.925
$ws.Cells.Item(20, 29).Value = -1
$ws.Cells.Item(20, 30).Value = 1.025

# Row 21
$ws.Cells.Item(21, 2).Value = 6799252
$ws.Cells.Item(21, 5).Value = 'AC Horsens'
$ws.Cells.Item(21, 6).Value = 'HB Kge'
$ws.Cells.Item(21, 7).Value = 2
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 2
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 'H'
$ws.Cells.Item(21, 12).Value = 1.7
$ws.Cells.Item(21, 13).Value = 3.6
$ws.Cells.Item(21, 14).Value = 4.75
$ws.Cells.Item(21, 15).Value = 1.727
$ws.Cells.Item(21, 16).Value = 4
$ws.Cells.Item(21, 17).Value = 4
$ws.Cells.Item(21, 18).Value = -0.75
$ws.Cells.Item(21, 19).Value = 2
$ws.Cells.Item(21, 20).Value = 1.85
$ws.Cells.Item(21, 21).Value = 3
$ws.Cells.Item(21, 22).Value = 1.95
$ws.Cells.Item(21, 23).Value = 1.9
$ws.Cells.Item(21, 24).Value = 0.7270000000000001
$ws.Cells.Item(21, 25).Value = -1
$ws.Cells.Item(21, 26).Value = -1
$ws.Cells.Item(21, 27).Value = 1
$ws.Cells.Item(21, 28).Value = -1
$ws.Cells.Item(21, 29).Value = -1
$ws.Cells.Item(21, 30).Value = 0.8999999999999999

# Row 40
$ws.Cells.Item(40, 2).Value = 6798552
$ws.Cells.Item(40, 5).Value = 'Vendsyssel FF'
$ws.Cells.Item(40, 6).Value = 'Kolding IF'
$ws.Cells.Item(40, 7).Value = 2
$ws.Cells.Item(40, 8).Value = 1
$ws.Cells.Item(40, 9).Value = 1
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 'H'
$ws.Cells.Item(40, 12).Value = 2.2
$ws.Cells.Item(40, 13).Value = 3.5
$ws.Cells.Item(40, 14).Value = 2.8
$ws.Cells.Item(40, 15).Value = 2.7
$ws.Cells.Item(40, 16).Value = 3.5
$ws.Cells.Item(40, 17).Value = 2.5
$ws.Cells.Item(40, 18).Value = 0
$ws.Cells.Item(40, 19).Value = 2
$ws.Cells.Item(40, 20).Value = 1.85
$ws.Cells.Item(40, 21).Value = 2.5
$ws.Cells.Item(40, 22).Value = 1.875
$ws.Cells.Item(40, 23).Value = 1.975
$ws.Cells.Item(40, 24).Value = 1.7
$ws.Cells.Item(40, 25).Value = -1
$ws.Cells.Item(40, 26).Value = -1
$ws.Cells.Item(40, 27).Value = 1
$ws.Cells.Item(40, 28).Value = -1
$ws.Cells.Item(40, 29).Value = 0.875
$ws.Cells.Item(40, 30).Value = -1

# Row 42
$ws.Cells.Item(42, 2).Value = 6799263
$ws.Cells.Item(42, 5).Value = 'FC Helsingor'
$ws.Cells.Item(42, 6).Value = 'Hillerd'
$ws.Cells.Item(42, 7).Value = 1
$ws.Cells.Item(42, 8).Value = 6
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 2
$ws.Cells.Item(42, 11).Value = 'A'
$ws.Cells.Item(42, 12).Value = 1.75
$ws.Cells.Item(42, 13).Value = 3.8
$ws.Cells.Item(42, 14).Value = 4
$ws.Cells.Item(42, 15).Value = 1.75
$ws.Cells.Item(42, 16).Value = 4
$ws.Cells.Item(42, 17).Value = 4.333
$ws.Cells.Item(42, 18).Value = -0.75
$ws.Cells.Item(42, 19).Value = 2
$ws.Cells.Item(42, 20).Value = 1.85
$ws.Cells.Item(42, 21).Value = 3
$ws.Cells.Item(42, 22).Value = 1.925
$ws.Cells.Item(42, 23).Value = 1.925
$ws.Cells.Item(42, 24).Value = -1
$ws.Cells.Item(42, 25).Value = -1
$ws.Cells.Item(42, 26).Value = 3.333
$ws.Cells.Item(42, 27).Value = -1
$ws.Cells.Item(42, 28).Value = 0.8500000000000001
$ws.Cells.Item(42, 29).Value = 0.925
$ws.Cells.Item(42, 30).Value = -1

# Row 98
$ws.Cells.Item(98, 2).Value = 6800819
$ws.Cells.Item(98, 5).Value = 'B93 Copenhagen'
$ws.Cells.Item(98, 6).Value = 'FC Helsingor'
$ws.Cells.Item(98, 7).Value = 2
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 'H'
$ws.Cells.Item(98, 12).Value = 2.5
$ws.Cells.Item(98, 13).Value = 3.4
$ws.Cells.Item(98, 14).Value = 2.55
$ws.Cells.Item(98, 15).Value = 2.45
$ws.Cells.Item(98, 16).Value = 3.5
$ws.Cells.Item(98, 17).Value = 2.8
$ws.Cells.Item(98, 18).Value = 0
$ws.Cells.Item(98, 19).Value = 1.8
$ws.Cells.Item(98, 20).Value = 2.05
$ws.Cells.Item(98, 21).Value = 2.75
$ws.Cells.Item(98, 22).Value = 1.9
$ws.Cells.Item(98, 23).Value = 1.95
$ws.Cells.Item(98, 24).Value = 1.45
$ws.Cells.Item(98, 25).Value = -1
$ws.Cells.Item(98, 26).Value = -1
$ws.Cells.Item(98, 27).Value = 0.8
$ws.Cells.Item(98, 28).Value = -1
$ws.Cells.Item(98, 29).Value = -1
$ws.Cells.Item(98, 30).Value = 0.95

# Row 99
$ws.Cells.Item(99, 2).Value = 6799296
$ws.Cells.Item(99, 5).Value = 'AC Horsens'
$ws.Cells.Item(99, 6).Value = 'FC Fredericia'
$ws.Cells.Item(99, 7).Value = 1
$ws.Cells.Item(99, 8).Value = 1
$ws.Cells.Item(99, 9).Value = 1
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 'D'
$ws.Cells.Item(99, 12).Value = 2.75
$ws.Cells.Item(99, 13).Value = 3.4
$ws.Cells.Item(99, 14).Value = 2.3
$ws.Cells.Item(99, 15).Value = 2.625
$ws.Cells.Item(99, 16).Value = 3.6
$ws.Cells.Item(99, 17).Value = 2.5
$ws.Cells.Item(99, 18).Value = 0
$ws.Cells.Item(99, 19).Value = 2
$ws.Cells.Item(99, 20).Value = 1.85
$ws.Cells.Item(99, 21).Value = 3
$ws.Cells.Item(99, 22).Value = 1.95
$ws.Cells.Item(99, 23).Value = 1.9
$ws.Cells.Item(99, 24).Value = -1
$ws.Cells.Item(99, 25).Value = 2.6
$ws.Cells.Item(99, 26).Value = -1
$ws.Cells.Item(99, 27).Value = 0
$ws.Cells.Item(99, 28).Value = 0
$ws.Cells.Item(99, 29).Value = -1
$ws.Cells.Item(99, 30).Value = 0.8999999999999999

# Row 114
$ws.Cells.Item(114, 2).Value = 6799302
$ws.Cells.Item(114, 5).Value = 'Naestved'
$ws.Cells.Item(114, 6).Value = 'FC Helsingor'
$ws.Cells.Item(114, 7).Value = 2
$ws.Cells.Item(114, 8).Value = 2
$ws.Cells.Item(114, 9).Value = 0
$ws.Cells.Item(114, 10).Value = 1
$ws.Cells.Item(114, 11).Value = 'D'
$ws.Cells.Item(114, 12).Value = 1.85
$ws.Cells.Item(114, 13).Value = 3.75
$ws.Cells.Item(114, 14).Value = 3.75
$ws.Cells.Item(114, 15).Value = 2.1
$ws.Cells.Item(114, 16).Value = 3.5
$ws.Cells.Item(114, 17).Value = 3.3
$ws.Cells.Item(114, 18).Value = -0.25
$ws.Cells.Item(114, 19).Value = 1.825
$ws.Cells.Item(114, 20).Value = 2.025
$ws.Cells.Item(114, 21).Value = 2.5
$ws.Cells.Item(114, 22).Value = 1.825
$ws.Cells.Item(114, 23).Value = 2.025
$ws.Cells.Item(114, 24).Value = -1
$ws.Cells.Item(114, 25).Value = 2.5
$ws.Cells.Item(114, 26).Value = -1
$ws.Cells.Item(114, 27).Value = -0.5
$ws.Cells.Item(114, 28).Value = 0.5125
$ws.Cells.Item(114, 29).Value = 0.825
$ws.Cells.Item(114, 30).Value = -1

# Row 115
$ws.Cells.Item(115, 2).Value = 7554250
$ws.Cells.Item(115, 5).Value = 'Kolding IF'
$ws.Cells.Item(115, 6).Value = 'B93 Copenhagen'
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 1
$ws.Cells.Item(115, 9).Value = 0
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 11).Value = 'A'
$ws.Cells.Item(115, 12).Value = 1.571
$ws.Cells.Item(115, 13).Value = 4.5
$ws.Cells.Item(115, 14).Value = 5
$ws.Cells.Item(115, 15).Value = 1.533
$ws.Cells.Item(115, 16).Value = 4.5
$ws.Cells.Item(115, 17).Value = 5.5
$ws.Cells.Item(115, 18).Value = -1
$ws.Cells.Item(115, 19).Value = 1.875
$ws.Cells.Item(115, 20).Value = 1.975
$ws.Cells.Item(115, 21).Value = 2.75
$ws.Cells.Item(115, 22).Value = 1.975
$ws.Cells.Item(115, 23).Value = 1.875
$ws.Cells.Item(115, 24).Value = -1
$ws.Cells.Item(115, 25).Value = -1
$ws.Cells.Item(115, 26).Value = 4.5
$ws.Cells.Item(115, 27).Value = -1
$ws.Cells.Item(115, 28).Value = 0.9750000000000001
$ws.Cells.Item(115, 29).Value = -1
$ws.Cells.Item(115, 30).Value = 0.875

# Row 116
$ws.Cells.Item(116, 2).Value = 6799307
$ws.Cells.Item(116, 5).Value = 'Vendsyssel FF'
$ws.Cells.Item(116, 6).Value = 'Hillerd'
$ws.Cells.Item(116, 7).Value = 4
$ws.Cells.Item(116, 8).Value = 4
$ws.Cells.Item(116, 9).Value = 1
$ws.Cells.Item(116, 10).Value = 2
$ws.Cells.Item(116, 11).Value = 'D'
$ws.Cells.Item(116, 12).Value = 1.75
$ws.Cells.Item(116, 13).Value = 3.6
$ws.Cells.Item(116, 14).Value = 4.333
$ws.Cells.Item(116, 15).Value = 2.4
$ws.Cells.Item(116, 16).Value = 3.2
$ws.Cells.Item(116, 17).Value = 3
$ws.Cells.Item(116, 18).Value = -0.25
$ws.Cells.Item(116, 19).Value = 2.1
$ws.Cells.Item(116, 20).Value = 1.775
$ws.Cells.Item(116, 21).Value = 2.25
$ws.Cells.Item(116, 22).Value = 1.85
$ws.Cells.Item(116, 23).Value = 2
$ws.Cells.Item(116, 24).Value = -1
$ws.Cells.Item(116, 25).Value = 2.2
$ws.Cells.Item(116, 26).Value = -1
$ws.Cells.Item(116, 27).Value = -0.5
$ws.Cells.Item(116, 28).Value = 0.3875
$ws.Cells.Item(116, 29).Value = 0.8500000000000001
$ws.Cells.Item(116, 30).Value = -1

# Row 118
$ws.Cells.Item(118, 2).Value = 6798562
$ws.Cells.Item(118, 5).Value = 'AC Horsens'
$ws.Cells.Item(118, 6).Value = 'Kolding IF'
$ws.Cells.Item(118, 7).Value = 1
$ws.Cells.Item(118, 8).Value = 2
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 11).Value = 'A'
$ws.Cells.Item(118, 12).Value = 2.4
$ws.Cells.Item(118, 13).Value = 3.5
$ws.Cells.Item(118, 14).Value = 2.4
$ws.Cells.Item(118, 15).Value = 3.4
$ws.Cells.Item(118, 16).Value = 3.3
$ws.Cells.Item(118, 17).Value = 2.15
$ws.Cells.Item(118, 18).Value = 0.25
$ws.Cells.Item(118, 19).Value = 1.975
$ws.Cells.Item(118, 20).Value = 1.875
$ws.Cells.Item(118, 21).Value = 2.25
$ws.Cells.Item(118, 22).Value = 1.975
$ws.Cells.Item(118, 23).Value = 1.875
$ws.Cells.Item(118, 24).Value = -1
$ws.Cells.Item(118, 25).Value = -1
$ws.Cells.Item(118, 26).Value = 1.15
$ws.Cells.Item(118, 27).Value = -1
$ws.Cells.Item(118, 28).Value = 0.875
$ws.Cells.Item(118, 29).Value = 0.9750000000000001
$ws.Cells.Item(118, 30).Value = -1

# Row 188
$ws.Cells.Item(188, 2).Value = 8231375
$ws.Cells.Item(188, 5).Value = 'Vendsyssel FF'
$ws.Cells.Item(188, 6).Value = 'FC Fredericia'
$ws.Cells.Item(188, 7).Value = 1
$ws.Cells.Item(188, 8).Value = 4
$ws.Cells.Item(188, 9).Value = 0
$ws.Cells.Item(188, 10).Value = 1
$ws.Cells.Item(188, 11).Value = 'A'
$ws.Cells.Item(188, 12).Value = 1.9
$ws.Cells.Item(188, 13).Value = 3.8
$ws.Cells.Item(188, 14).Value = 3.25
$ws.Cells.Item(188, 15).Value = 2.3
$ws.Cells.Item(188, 16).Value = 3.8
$ws.Cells.Item(188, 17).Value = 2.7
$ws.Cells.Item(188, 18).Value = -0.25
$ws.Cells.Item(188, 19).Value = 2.05
$ws.Cells.Item(188, 20).Value = 1.8
$ws.Cells.Item(188, 21).Value = 3.25
$ws.Cells.Item(188, 22).Value = 1.95
$ws.Cells.Item(188, 23).Value = 1.9
$ws.Cells.Item(188, 24).Value = -1
$ws.Cells.Item(188, 25).Value = -1
$ws.Cells.Item(188, 26).Value = 1.7
$ws.Cells.Item(188, 27).Value = -1
$ws.Cells.Item(188, 28).Value = 0.8
$ws.Cells.Item(188, 29).Value = 0.95
$ws.Cells.Item(188, 30).Value = -1

# Row 189
$ws.Cells.Item(189, 2).Value = 8231374
$ws.Cells.Item(189, 5).Value = 'Kolding IF'
$ws.Cells.Item(189, 6).Value = 'AaB'
$ws.Cells.Item(189, 7).Value = 2
$ws.Cells.Item(189, 8).Value = 3
$ws.Cells.Item(189, 9).Value = 1
$ws.Cells.Item(189, 10).Value = 1
$ws.Cells.Item(189, 11).Value = 'A'
$ws.Cells.Item(189, 12).Value = 2.5
$ws.Cells.Item(189, 13).Value = 3.5
$ws.Cells.Item(189, 14).Value = 2.6
$ws.Cells.Item(189, 15).Value = 2.1
$ws.Cells.Item(189, 16).Value = 3.6
$ws.Cells.Item(189, 17).Value = 3.1
$ws.Cells.Item(189, 18).Value = -0.25
$ws.Cells.Item(189, 19).Value = 1.875
$ws.Cells.Item(189, 20).Value = 1.975
$ws.Cells.Item(189, 21).Value = 2.75
$ws.Cells.Item(189, 22).Value = 1.85
$ws.Cells.Item(189, 23).Value = 2
$ws.Cells.Item(189, 24).Value = -1
$ws.Cells.Item(189, 25).Value = -1
$ws.Cells.Item(189, 26).Value = 2.1
$ws.Cells.Item(189, 27).Value = -1
$ws.Cells.Item(189, 28).Value = 0.9750000000000001
$ws.Cells.Item(189, 29).Value = 0.8500000000000001
$ws.Cells.Item(189, 30).Value = -1

# Row 190
$ws.Cells.Item(190, 2).Value = 8231242
$ws.Cells.Item(190, 5).Value = 'Hobro IK'
$ws.Cells.Item(190, 6).Value = 'Sonderjyske'
$ws.Cells.Item(190, 7).Value = 2
$ws.Cells.Item(190, 8).Value = 2
$ws.Cells.Item(190, 9).Value = 0
$ws.Cells.Item(190, 10).Value = 1
$ws.Cells.Item(190, 11).Value = 'D'
$ws.Cells.Item(190, 12).Value = 4
$ws.Cells.Item(190, 13).Value = 3.8
$ws.Cells.Item(190, 14).Value = 1.727
$ws.Cells.Item(190, 15).Value = 4
$ws.Cells.Item(190, 16).Value = 3.75
$ws.Cells.Item(190, 17).Value = 1.8
$ws.Cells.Item(190, 18).Value = 0.5
$ws.Cells.Item(190, 19).Value = 2
$ws.Cells.Item(190, 20).Value = 1.85
$ws.Cells.Item(190, 21).Value = 3
$ws.Cells.Item(190, 22).Value = 1.925
$ws.Cells.Item(190, 23).Value = 1.925
$ws.Cells.Item(190, 24).Value = -1
$ws.Cells.Item(190, 25).Value = 2.75
$ws.Cells.Item(190, 26).Value = -1
$ws.Cells.Item(190, 27).Value = 1
$ws.Cells.Item(190, 28).Value = -1
$ws.Cells.Item(190, 29).Value = 0.925
$ws.Cells.Item(190, 30).Value = -1
